$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9: title/link update
$ws.Range("D9").Value = "데이터 사이언스 (Data Science) 석사과정 강의목록"
$ws.Range("E9").Value = "https://blog.pabii.co.kr/msds-course-works/#utm_source=rss&utm_medium=rss&utm_campaign=msds-course-works"

# Row 26: title/link update
$ws.Range("D26").Value = "AI를 활용한 고용률 예측 모델 개발기(1)"
$ws.Range("E26").Value = "https://blog.est.ai/2021/03/employment-rate/"

# Row 32: title/link update
$ws.Range("D32").Value = "Autoencoder 를 이용한 차원 축소 (latent representation)"
$ws.Range("E32").Value = "https://dodonam.tistory.com/301"

# Row 39: title/link update
$ws.Range("D39").Value = "Probability concepts explained: Maximum likelihood estimation"
$ws.Range("E39").Value = "https://a292run.tistory.com/entry/Probability-concepts-explained-Maximum-likelihood-estimation-1"

# Row 51: title update only
$ws.Range("D51").Value = "블루스킨(Blue skin v1.1)을 소개합니다"
